$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.343.77"
$ws.Range("E2").Value = "  -2.14%  "

$ws.Range("D3").Value = "1.865.17"
$ws.Range("E3").Value = "  -1.95%  "

$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").Value = "319.84"
$ws.Range("E5").Value = "  -1.42%  "

$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("D7").Value = "0.4409"
$ws.Range("E7").Value = "  -4.09%  "

$ws.Range("D8").Value = "0.3723"
$ws.Range("E8").Value = "  -2.33%  "

$ws.Range("D9").Value = "0.07558"
$ws.Range("E9").Value = "  -1.89%  "

$ws.Range("D10").Value = "0.9399"
$ws.Range("E10").Value = "  -3.37%  "

$ws.Range("D11").Value = "21.33"
$ws.Range("E11").Value = "  -2.65%  "

$ws.Range("D12").Value = "1.847.10"
$ws.Range("E12").Value = "  -2.93%  "

$ws.Range("D13").Value = "6.744"
$ws.Range("E13").Value = "  -2.57%  "

$ws.Range("D14").Value = "5.483"
$ws.Range("E14").Value = "  -2.86%  "

$ws.Range("D15").Value = "0.06865"
$ws.Range("E15").Value = "  -3.01%  "

$ws.Range("E16").Value = "  -0.16%  "

$ws.Range("D17").Value = "82.18"
$ws.Range("E17").Value = "  -1.88%  "

$ws.Range("D18").Value = "0.000009116"
$ws.Range("E18").Value = "  -3.89%  "

$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("E20").Value = "  -3.61%  "

$ws.Range("D21").Value = "28.332.88"
$ws.Range("E21").Value = "  -2.05%  "

$ws.Range("D22").Value = "5.154"
$ws.Range("E22").Value = "  -2.62%  "

$ws.Range("D23").Value = "10.77"
$ws.Range("E23").Value = "  -0.79%  "

$ws.Range("D24").Value = "2.110.09"
$ws.Range("E24").Value = "  -1.65%  "

$ws.Range("D25").Value = "2.041"
$ws.Range("E25").Value = "  -2.76%  "

$ws.Range("D26").Value = "154.77"
$ws.Range("E26").Value = "  -2.02%  "

$ws.Range("D27").Value = "18.47"
$ws.Range("E27").Value = "  -3.17%  "

$ws.Range("D28").Value = "5.373"
$ws.Range("E28").Value = "  -4.34%  "

$ws.Range("D29").Value = "114.69"
$ws.Range("E29").Value = "  -2.44%  "

$ws.Range("D30").Value = "1.734"
$ws.Range("E30").Value = "  -6.07%  "

$ws.Range("D31").Value = "0.09082"
$ws.Range("E31").Value = "  -1.87%  "

$ws.Range("D32").Value = "0.8079"
$ws.Range("E32").Value = "  -5.80%  "

$ws.Range("D33").Value = "4.869"
$ws.Range("E33").Value = "  -4.32%  "

$ws.Range("D34").Value = "1.180"
$ws.Range("E34").Value = "  -4.66%  "

$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("E36").Value = "  -0.17%  "

$ws.Range("D37").Value = "1.127"
$ws.Range("E37").Value = "  -1.13%  "

$ws.Range("D38").Value = "0.05489"
$ws.Range("E38").Value = "  -3.32%  "

$ws.Range("D39").Value = "3.014"
$ws.Range("E39").Value = "  +9.37%  "

$ws.Range("D40").Value = "0.01957"
$ws.Range("E40").Value = "  -3.62%  "

$ws.Range("D41").Value = "7.189"
$ws.Range("E41").Value = "  -2.62%  "

$ws.Range("D42").Value = "0.5267"
$ws.Range("E42").Value = "  -3.78%  "

$ws.Range("E43").Value = "  -4.12%  "

$ws.Range("D44").Value = "8.879"
$ws.Range("E44").Value = "  -4.36%  "

$ws.Range("D45").Value = "0.06782"
$ws.Range("E45").Value = "  -0.71%  "

$ws.Range("D46").Value = "2.061"
$ws.Range("E46").Value = "  -0.41%  "

$ws.Range("D47").Value = "0.4903"
$ws.Range("E47").Value = "  -4.83%  "

$ws.Range("D48").Value = "0.000002540"
$ws.Range("E48").Value = "  -1.05%  "

$ws.Range("D49").Value = "10.60"
$ws.Range("E49").Value = "  -4.75%  "

$ws.Range("D50").Value = "107.80"
$ws.Range("E50").Value = "  -2.00%  "

$ws.Range("D51").Value = "1.690"
$ws.Range("E51").Value = "  -4.37%  "
